$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at the beginning (A and B), shifting old A:C to C:E
$ws.Range("A1:B1").EntireColumn.Insert()

# Set header labels for the two new columns
$ws.Range("A1").Value = "vendor"
$ws.Range("B1").Value = "doc. number"

# Set column widths to match target layout (values chosen so the runtime's
# internal char->stored-width rounding lands as close as possible to the
# target stored widths of 14.85546875 / 11.5703125 / 16.7109375 / 16)
$ws.Columns("B").ColumnWidth = 14
$ws.Columns("C").ColumnWidth = 10.666666666666666
$ws.Columns("D").ColumnWidth = 15.833333333333334
$ws.Columns("E").ColumnWidth = 15.166666666666666

# Fill in "vendor" values (column A)
$ws.Range("A2").Value = 8946516
$ws.Range("A3").Value = 654513
$ws.Range("A4").Value = 6543210
$ws.Range("A5").Value = 6846321
$ws.Range("A6").Value = 5644668
$ws.Range("A7").Value = 5603479.2000000002
$ws.Range("A8").Value = 5562290.4000000004
$ws.Range("A9").Value = 5521101.5999999996
$ws.Range("A10").Value = 5479912.7999999998

# Fill in "doc. number" values (column B)
$ws.Range("B2").Value = 74000056
$ws.Range("B3").Value = 74000085
$ws.Range("B4").Value = 74000114
$ws.Range("B5").Value = 74000143
$ws.Range("B6").Value = 74000172
$ws.Range("B7").Value = 74000201
$ws.Range("B8").Value = 74000230
$ws.Range("B9").Value = 74000259
$ws.Range("B10").Value = 74000288

# Update selection to match target
[void]$ws.Range("G9").Select()
